# Add a new "unknow-message" row (tag + long apology response) right after the
# existing unknow-message row, and give the "law-place" row an updated answer
# that now starts with an introductory sentence. Net effect: one row is moved
# from the bottom of the sheet up to row 8, the law-place response text is
# replaced, and a new blank formatted row is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the last row (unknow-message / long apology text) and move
#        it up to become the new row 8, right after the current row 7
#        (unknow-message / "ไม่เข้าใจคำถาม"). ---
$ws.Rows.Item(31).Copy()
$ws.Rows.Item(8).Insert()

# The row that used to be row 31 is now pushed down to row 32 as a duplicate;
# remove that duplicate now that its content lives at row 8.
$ws.Rows.Item(32).Delete()

# Row height isn't carried over by Insert, so restore it explicitly to match
# the source row's auto-computed wrap height.
$ws.Rows.Item(8).RowHeight = 242.25

# --- 2. Update the "law-place" row (now row 30, after the insert above
#        shifted it down by one) with the new, longer response text that
#        adds an introductory sentence before the existing three bullet
#        points. ---
$ws.Cells.Item(30, 2).Value = "สามารถทำการยื่นเสียภาษีได้ที่สถานที่ต่อไปนี้เลยค่ะ" + "`n" + "1.สำนักงานสรรพากรทุกสาขาทุกเเห่ง " + "`n" + "2.ไปรษณีย์ เเบบลงทะเบียน " + "`n" + "3.ช่องทางออนไลน์ ผ่านเว็บไซต์ของกรมสรรพากร"
$ws.Rows.Item(30).RowHeight = 63.75

# --- 3. Append a new, empty but formatted row 32 (same look as the other
#        answer cells in column B) so the sheet ends with a blank row. ---
$ws.Range("B9").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 4. Restore the original selection/view state. ---
$ws.Application.Goto($ws.Range("C26"), $false)
$ws.Cells.Item(26, 3).Select()
